$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 90, shifting existing rows 90-99 down to 91-100.
$ws.Rows.Item(90).Insert()

# Populate the new weekly record in row 90.
$ws.Cells.Item(90, 1).Value = 10
$ws.Cells.Item(90, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(90, 3).Value = "La Araucanía"
$ws.Cells.Item(90, 4).Value = 45041
$ws.Cells.Item(90, 5).Value = 9
$ws.Cells.Item(90, 6).Value = 300000001
$ws.Cells.Item(90, 7).Value = "Rabanito"
$ws.Cells.Item(90, 8).Value = "Sin especificar"
$ws.Cells.Item(90, 9).Value = "Primera"
$ws.Cells.Item(90, 10).Value = 50
$ws.Cells.Item(90, 11).Value = 7000
$ws.Cells.Item(90, 12).Value = 8000
$ws.Cells.Item(90, 13).Value = 7500
$ws.Cells.Item(90, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(90, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(90, 16).Value = 625
$ws.Cells.Item(90, 17).Value = 12
$ws.Cells.Item(90, 18).Value = "Hortaliza"
